$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.400.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3695"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3278"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07406"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.834"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.804"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.553.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06698"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.07"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.329"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.389.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.313"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.577"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.43%  "

$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.942"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.735.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.044"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.973"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.965"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.655"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08240"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02387"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.295"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06305"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2185"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.227"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6109"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "

$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5947"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.743"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.011"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.180"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07155"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.78%  "
